$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "14 Nov 2025, 09:17 AM"

# --- Top Gainers sheet: rows 21-76 shift up by one (oldest dropped off, new entrant appended) ---
$gainers = $wb.Worksheets.Item("Top Gainers")
$gainers.Range("B21").Value = "SURAJEST"
$gainers.Range("C21").Value = 3.4309
$gainers.Range("D21").Value = -0.8438
$gainers.Range("E21").Value = -0.5011
$gainers.Range("B22").Value = "SRM"
$gainers.Range("C22").Value = 3.4239
$gainers.Range("D22").Value = 4.5145
$gainers.Range("E22").Value = 4.9444
$gainers.Range("B23").Value = "RUBICON"
$gainers.Range("C23").Value = 3.301
$gainers.Range("D23").Value = 7.7726
$gainers.Range("E23").Value = 9.3505
$gainers.Range("B24").Value = "KIOCL"
$gainers.Range("C24").Value = 3.2004
$gainers.Range("D24").Value = -0.757
$gainers.Range("E24").Value = -7.1735
$gainers.Range("B25").Value = "RNBDENIMS"
$gainers.Range("C25").Value = 2.9689
$gainers.Range("D25").Value = -12.2615
$gainers.Range("E25").Value = -5.4969
$gainers.Range("B26").Value = "HONDAPOWER"
$gainers.Range("C26").Value = 2.7729
$gainers.Range("D26").Value = 5.3636
$gainers.Range("E26").Value = 1.4238
$gainers.Range("B27").Value = "UNIVCABLES"
$gainers.Range("C27").Value = 2.7551
$gainers.Range("D27").Value = 11.1111
$gainers.Range("E27").Value = 9.014099999999999
$gainers.Range("B28").Value = "PRABHA"
$gainers.Range("C28").Value = 2.7113
$gainers.Range("D28").Value = 3.9652
$gainers.Range("E28").Value = 0.0564
$gainers.Range("B29").Value = "DELTACORP"
$gainers.Range("C29").Value = 2.5891
$gainers.Range("D29").Value = -0.0678
$gainers.Range("E29").Value = -5.2212
$gainers.Range("B30").Value = "MATRIMONY"
$gainers.Range("C30").Value = 2.3672
$gainers.Range("D30").Value = 0.2826
$gainers.Range("E30").Value = 0.9118000000000001
$gainers.Range("B31").Value = "GMRAIRPORT"
$gainers.Range("C31").Value = 2.262
$gainers.Range("D31").Value = 2.0803
$gainers.Range("E31").Value = 3.9161
$gainers.Range("B32").Value = "RAMASTEEL"
$gainers.Range("C32").Value = 2.1556
$gainers.Range("D32").Value = 3.4156
$gainers.Range("E32").Value = 6.5494
$gainers.Range("B33").Value = "SANDHAR"
$gainers.Range("C33").Value = 2.1423
$gainers.Range("D33").Value = 3.6032
$gainers.Range("E33").Value = 10.33
$gainers.Range("B34").Value = "RUPA"
$gainers.Range("C34").Value = 2.0616
$gainers.Range("D34").Value = 4.0585
$gainers.Range("E34").Value = 2.2073
$gainers.Range("B35").Value = "LGHL"
$gainers.Range("C35").Value = 2.0444
$gainers.Range("D35").Value = 10.4909
$gainers.Range("E35").Value = 5.6118
$gainers.Range("B36").Value = "SAURASHCEM"
$gainers.Range("C36").Value = 2.0202
$gainers.Range("D36").Value = 2.3963
$gainers.Range("E36").Value = -1.0392
$gainers.Range("B37").Value = "MOTHERSON"
$gainers.Range("C37").Value = 1.9976
$gainers.Range("D37").Value = 7.2351
$gainers.Range("E37").Value = 5.6172
$gainers.Range("B38").Value = "JINDALPHOT"
$gainers.Range("C38").Value = 1.9875
$gainers.Range("D38").Value = 2.9057
$gainers.Range("E38").Value = 5.5001
$gainers.Range("B39").Value = "RAJRILTD"
$gainers.Range("C39").Value = 1.9802
$gainers.Range("D39").Value = 7.3364
$gainers.Range("E39").Value = 5.6627
$gainers.Range("B40").Value = "MONTECARLO"
$gainers.Range("C40").Value = 1.9442
$gainers.Range("D40").Value = 10.4863
$gainers.Range("E40").Value = 11.826
$gainers.Range("B41").Value = "AEROENTER"
$gainers.Range("C41").Value = 1.8601
$gainers.Range("D41").Value = 2.0748
$gainers.Range("E41").Value = 5.6985
$gainers.Range("B42").Value = "EBGNG"
$gainers.Range("C42").Value = 1.8574
$gainers.Range("D42").Value = 1.5405
$gainers.Range("E42").Value = 2.3287
$gainers.Range("B43").Value = "NUVOCO"
$gainers.Range("C43").Value = 1.8474
$gainers.Range("D43").Value = -1.0092
$gainers.Range("E43").Value = -8.9232
$gainers.Range("B44").Value = "SANOFICONR"
$gainers.Range("C44").Value = 1.8095
$gainers.Range("D44").Value = 9.313700000000001
$gainers.Range("E44").Value = 4.7991
$gainers.Range("B45").Value = "IFGLEXPOR"
$gainers.Range("C45").Value = 1.7977
$gainers.Range("D45").Value = -0.1572
$gainers.Range("E45").Value = -5.0121
$gainers.Range("B46").Value = "PDMJEPAPER"
$gainers.Range("C46").Value = 1.7332
$gainers.Range("D46").Value = -5.5366
$gainers.Range("E46").Value = -6.6471
$gainers.Range("B47").Value = "HINDOILEXP"
$gainers.Range("C47").Value = 1.6061
$gainers.Range("D47").Value = 0.1946
$gainers.Range("E47").Value = -2.4438
$gainers.Range("B48").Value = "DIVGIITTS"
$gainers.Range("C48").Value = 1.5922
$gainers.Range("D48").Value = 6.48
$gainers.Range("E48").Value = 5.5946
$gainers.Range("B49").Value = "INFIBEAM"
$gainers.Range("C49").Value = 1.5625
$gainers.Range("D49").Value = 6.7323
$gainers.Range("E49").Value = 0
$gainers.Range("B50").Value = "BLUSPRING"
$gainers.Range("C50").Value = 1.5196
$gainers.Range("D50").Value = 4.5502
$gainers.Range("E50").Value = 0.2902
$gainers.Range("B51").Value = "MANAPPURAM"
$gainers.Range("C51").Value = 1.5154
$gainers.Range("D51").Value = 2.7157
$gainers.Range("E51").Value = 3.154
$gainers.Range("B52").Value = "SKYGOLD"
$gainers.Range("C52").Value = 1.5081
$gainers.Range("D52").Value = -0.5682
$gainers.Range("E52").Value = -1.644
$gainers.Range("B53").Value = "SANATHAN"
$gainers.Range("C53").Value = 1.5072
$gainers.Range("D53").Value = -1.6031
$gainers.Range("E53").Value = 1.0208
$gainers.Range("B54").Value = "SANGAMIND"
$gainers.Range("C54").Value = 1.5047
$gainers.Range("D54").Value = -4.8468
$gainers.Range("E54").Value = -0.4549
$gainers.Range("B55").Value = "FINKURVE"
$gainers.Range("C55").Value = 1.5009
$gainers.Range("D55").Value = 11.7543
$gainers.Range("E55").Value = 5.3975
$gainers.Range("B56").Value = "TVSELECT"
$gainers.Range("C56").Value = 1.5009
$gainers.Range("D56").Value = 12.6287
$gainers.Range("E56").Value = 9.7895
$gainers.Range("B57").Value = "INNOVACAP"
$gainers.Range("C57").Value = 1.4998
$gainers.Range("D57").Value = -3.4463
$gainers.Range("E57").Value = -10.835
$gainers.Range("B58").Value = "RMDRIP"
$gainers.Range("C58").Value = 1.4889
$gainers.Range("D58").Value = -2.4107
$gainers.Range("E58").Value = -2.3868
$gainers.Range("B59").Value = "GLAND"
$gainers.Range("C59").Value = 1.4651
$gainers.Range("D59").Value = -0.3757
$gainers.Range("E59").Value = -3.7794
$gainers.Range("B60").Value = "SPIC"
$gainers.Range("C60").Value = 1.4347
$gainers.Range("D60").Value = 1.4802
$gainers.Range("E60").Value = -0.3852
$gainers.Range("B61").Value = "OSWALPUMPS"
$gainers.Range("C61").Value = 1.4151
$gainers.Range("D61").Value = -8.5755
$gainers.Range("E61").Value = -11.1448
$gainers.Range("B62").Value = "MARKSANS"
$gainers.Range("C62").Value = 1.407
$gainers.Range("D62").Value = 4.4001
$gainers.Range("E62").Value = 4.8989
$gainers.Range("B63").Value = "CPPLUS"
$gainers.Range("C63").Value = 1.392
$gainers.Range("D63").Value = 8.6716
$gainers.Range("E63").Value = 26.9406
$gainers.Range("B64").Value = "BHAGCHEM"
$gainers.Range("C64").Value = 1.3856
$gainers.Range("D64").Value = 4.3954
$gainers.Range("E64").Value = 2.503
$gainers.Range("B65").Value = "TVSSCS"
$gainers.Range("C65").Value = 1.3814
$gainers.Range("D65").Value = -0.063
$gainers.Range("E65").Value = -1.6118
$gainers.Range("B66").Value = "SASKEN"
$gainers.Range("C66").Value = 1.3761
$gainers.Range("D66").Value = -3.6035
$gainers.Range("E66").Value = 4.5259
$gainers.Range("B67").Value = "AHLUCONT"
$gainers.Range("C67").Value = 1.344
$gainers.Range("D67").Value = 0.4497
$gainers.Range("E67").Value = -5.9965
$gainers.Range("B68").Value = "ANUP"
$gainers.Range("C68").Value = 1.3321
$gainers.Range("D68").Value = 7.2157
$gainers.Range("E68").Value = 7.4853
$gainers.Range("B69").Value = "ALKEM"
$gainers.Range("C69").Value = 1.3099
$gainers.Range("D69").Value = 2.3287
$gainers.Range("E69").Value = 5.3391
$gainers.Range("B70").Value = "KICL"
$gainers.Range("C70").Value = 1.2885
$gainers.Range("D70").Value = 1.7589
$gainers.Range("E70").Value = -5.9079
$gainers.Range("B71").Value = "SANSERA"
$gainers.Range("C71").Value = 1.2862
$gainers.Range("D71").Value = 3.9741
$gainers.Range("E71").Value = 4.0497
$gainers.Range("B72").Value = "UNIPARTS"
$gainers.Range("C72").Value = 1.2814
$gainers.Range("D72").Value = 1.9292
$gainers.Range("E72").Value = -0.9026999999999999
$gainers.Range("B73").Value = "TVSHLTD"
$gainers.Range("C73").Value = 1.2668
$gainers.Range("D73").Value = 0.6461
$gainers.Range("E73").Value = -0.4349
$gainers.Range("B74").Value = "VMM"
$gainers.Range("C74").Value = 1.2594
$gainers.Range("D74").Value = -2.2908
$gainers.Range("E74").Value = -3.3105
$gainers.Range("B75").Value = "AJAXENGG"
$gainers.Range("C75").Value = 1.2439
$gainers.Range("D75").Value = -2.0616
$gainers.Range("E75").Value = -5.1426
$gainers.Range("B76").Value = "EUROPRATIK"
$gainers.Range("C76").Value = 1.2326
$gainers.Range("D76").Value = 6.0137
$gainers.Range("E76").Value = 9.8056

# --- Top Losers sheet: rows 16-76 shift up by one (oldest dropped off, new entrant appended) ---
$losers = $wb.Worksheets.Item("Top Losers")
$losers.Range("B16").Value = "STALLION"
$losers.Range("C16").Value = -3.3266
$losers.Range("D16").Value = -3.8023
$losers.Range("E16").Value = -21.6147
$losers.Range("B17").Value = "GHCL"
$losers.Range("C17").Value = -3.3063
$losers.Range("D17").Value = 0.2587
$losers.Range("E17").Value = -3.1553
$losers.Range("B18").Value = "SCODATUBES"
$losers.Range("C18").Value = -3.044
$losers.Range("D18").Value = 4.2867
$losers.Range("E18").Value = 7.6548
$losers.Range("B19").Value = "KINGFA"
$losers.Range("C19").Value = -3.0097
$losers.Range("D19").Value = -4.8396
$losers.Range("E19").Value = -8.3925
$losers.Range("B20").Value = "SUNFLAG"
$losers.Range("C20").Value = -2.967
$losers.Range("D20").Value = -0.1344
$losers.Range("E20").Value = -2.9489
$losers.Range("B21").Value = "WALCHANNAG"
$losers.Range("C21").Value = -2.786
$losers.Range("D21").Value = -4.0318
$losers.Range("E21").Value = -7.5692
$losers.Range("B22").Value = "MIDHANI"
$losers.Range("C22").Value = -2.7451
$losers.Range("D22").Value = -0.7603
$losers.Range("E22").Value = -3.4143
$losers.Range("B23").Value = "NEWGEN"
$losers.Range("C23").Value = -2.7117
$losers.Range("D23").Value = -1.0285
$losers.Range("E23").Value = -3.8828
$losers.Range("B24").Value = "DLINKINDIA"
$losers.Range("C24").Value = -2.633
$losers.Range("D24").Value = 2.4887
$losers.Range("E24").Value = -0.4177
$losers.Range("B25").Value = "MAITHANALL"
$losers.Range("C25").Value = -2.5933
$losers.Range("D25").Value = -0.0649
$losers.Range("E25").Value = -2.3904
$losers.Range("B26").Value = "PETRONET"
$losers.Range("C26").Value = -2.5099
$losers.Range("D26").Value = -2.3698
$losers.Range("E26").Value = -3.3244
$losers.Range("B27").Value = "VOLTAS"
$losers.Range("C27").Value = -2.4609
$losers.Range("D27").Value = -1.3243
$losers.Range("E27").Value = -5.7531
$losers.Range("B28").Value = "HERANBA"
$losers.Range("C28").Value = -2.446
$losers.Range("D28").Value = -4.1696
$losers.Range("E28").Value = -9.5548
$losers.Range("B29").Value = "RHETAN"
$losers.Range("C29").Value = -2.3929
$losers.Range("D29").Value = 0.6058
$losers.Range("E29").Value = -0.2574
$losers.Range("B30").Value = "PARACABLES"
$losers.Range("C30").Value = -2.3618
$losers.Range("D30").Value = -3.86
$losers.Range("E30").Value = -7.3261
$losers.Range("B31").Value = "BERGEPAINT"
$losers.Range("C31").Value = -2.3187
$losers.Range("D31").Value = 6.2988
$losers.Range("E31").Value = 4.1321
$losers.Range("B32").Value = "BHARATSE"
$losers.Range("C32").Value = -2.2676
$losers.Range("D32").Value = -10.3907
$losers.Range("E32").Value = -11.7457
$losers.Range("B33").Value = "AMNPLST"
$losers.Range("C33").Value = -2.2005
$losers.Range("D33").Value = -3.9431
$losers.Range("E33").Value = -7.7533
$losers.Range("B34").Value = "MAXIND"
$losers.Range("C34").Value = -2.1667
$losers.Range("D34").Value = -2.8895
$losers.Range("E34").Value = -3.5928
$losers.Range("B35").Value = "INDIAGLYCO"
$losers.Range("C35").Value = -2.1354
$losers.Range("D35").Value = 1.9401
$losers.Range("E35").Value = 5.4424
$losers.Range("B36").Value = "ONGC"
$losers.Range("C36").Value = -2.1328
$losers.Range("D36").Value = -2.6566
$losers.Range("E36").Value = -3.865
$losers.Range("B37").Value = "SAFARI"
$losers.Range("C37").Value = -2.1165
$losers.Range("D37").Value = 7.6886
$losers.Range("E37").Value = 3.3352
$losers.Range("B38").Value = "ADVENZYMES"
$losers.Range("C38").Value = -2.0653
$losers.Range("D38").Value = 7.6997
$losers.Range("E38").Value = 8.5213
$losers.Range("B39").Value = "VEEDOL"
$losers.Range("C39").Value = -2.0603
$losers.Range("D39").Value = 0.0808
$losers.Range("E39").Value = -1.4372
$losers.Range("B40").Value = "NATIONALUM"
$losers.Range("C40").Value = -2.036
$losers.Range("D40").Value = 12.1288
$losers.Range("E40").Value = 12.4114
$losers.Range("B41").Value = "SALZERELEC"
$losers.Range("C41").Value = -1.9807
$losers.Range("D41").Value = 5.8292
$losers.Range("E41").Value = -9.5886
$losers.Range("B42").Value = "TITAGARH"
$losers.Range("C42").Value = -1.9699
$losers.Range("D42").Value = 0.3202
$losers.Range("E42").Value = -4.3581
$losers.Range("B43").Value = "PLATIND"
$losers.Range("C43").Value = -1.8692
$losers.Range("D43").Value = -0.8499
$losers.Range("E43").Value = -4.9774
$losers.Range("B44").Value = "TEMBO"
$losers.Range("C44").Value = -1.8516
$losers.Range("D44").Value = 5.1308
$losers.Range("E44").Value = 13.0166
$losers.Range("B45").Value = "SPMLINFRA"
$losers.Range("C45").Value = -1.8097
$losers.Range("D45").Value = -7.0431
$losers.Range("E45").Value = -10.9693
$losers.Range("B46").Value = "DBREALTY"
$losers.Range("C46").Value = -1.7875
$losers.Range("D46").Value = 2.923
$losers.Range("E46").Value = 0.7523
$losers.Range("B47").Value = "IVALUE"
$losers.Range("C47").Value = -1.7748
$losers.Range("D47").Value = 9.0724
$losers.Range("E47").Value = 12.2574
$losers.Range("B48").Value = "VIDHIING"
$losers.Range("C48").Value = -1.7274
$losers.Range("D48").Value = 3.9619
$losers.Range("E48").Value = 0.1627
$losers.Range("B49").Value = "KRN"
$losers.Range("C49").Value = -1.7125
$losers.Range("D49").Value = -0.5712
$losers.Range("E49").Value = 0.6161
$losers.Range("B50").Value = "MINDACORP"
$losers.Range("C50").Value = -1.7015
$losers.Range("D50").Value = 2.7365
$losers.Range("E50").Value = 5.5615
$losers.Range("B51").Value = "AKUMS"
$losers.Range("C51").Value = -1.6918
$losers.Range("D51").Value = -5.1505
$losers.Range("E51").Value = -6.6739
$losers.Range("B52").Value = "FAZE3Q"
$losers.Range("C52").Value = -1.6553
$losers.Range("D52").Value = 5.5602
$losers.Range("E52").Value = 3.4624
$losers.Range("B53").Value = "POLYPLEX"
$losers.Range("C53").Value = -1.6408
$losers.Range("D53").Value = -1.412
$losers.Range("E53").Value = -4.7161
$losers.Range("B54").Value = "ORKLAINDIA"
$losers.Range("C54").Value = -1.6169
$losers.Range("D54").Value = -7.4758
$losers.Range("E54").Value = "N/A"
$losers.Range("B55").Value = "NPST"
$losers.Range("C55").Value = -1.6071
$losers.Range("D55").Value = 1.6949
$losers.Range("E55").Value = -3.059
$losers.Range("B56").Value = "SYSTMTXC"
$losers.Range("C56").Value = -1.5893
$losers.Range("D56").Value = -5.269
$losers.Range("E56").Value = -6.0737
$losers.Range("B57").Value = "INDIGOPNTS"
$losers.Range("C57").Value = -1.5809
$losers.Range("D57").Value = 26.4814
$losers.Range("E57").Value = 27.0127
$losers.Range("B58").Value = "ITI"
$losers.Range("C58").Value = -1.548
$losers.Range("D58").Value = 1.0968
$losers.Range("E58").Value = -2.0634
$losers.Range("B59").Value = "BAYERCROP"
$losers.Range("C59").Value = -1.5436
$losers.Range("D59").Value = -3.0445
$losers.Range("E59").Value = -8.8544
$losers.Range("B60").Value = "SOLARA"
$losers.Range("C60").Value = -1.5158
$losers.Range("D60").Value = -0.5397999999999999
$losers.Range("E60").Value = -6.3021
$losers.Range("B61").Value = "VAIBHAVGBL"
$losers.Range("C61").Value = -1.5077
$losers.Range("D61").Value = -1.7115
$losers.Range("E61").Value = -2.9238
$losers.Range("B62").Value = "SPAL"
$losers.Range("C62").Value = -1.5007
$losers.Range("D62").Value = 8.687900000000001
$losers.Range("E62").Value = -1.0412
$losers.Range("B63").Value = "SPECTRUM"
$losers.Range("C63").Value = -1.5
$losers.Range("D63").Value = 7.5987
$losers.Range("E63").Value = -1.558
$losers.Range("B64").Value = "LUMAXIND"
$losers.Range("C64").Value = -1.4996
$losers.Range("D64").Value = 5.5961
$losers.Range("E64").Value = -2.8732
$losers.Range("B65").Value = "SURAKSHA"
$losers.Range("C65").Value = -1.4995
$losers.Range("D65").Value = -4.4424
$losers.Range("E65").Value = -3.9314
$losers.Range("B66").Value = "ARTEMISMED"
$losers.Range("C66").Value = -1.4976
$losers.Range("D66").Value = 9.561500000000001
$losers.Range("E66").Value = 8.401999999999999
$losers.Range("B67").Value = "RAJOOENG"
$losers.Range("C67").Value = -1.4486
$losers.Range("D67").Value = -6.2148
$losers.Range("E67").Value = -6.9736
$losers.Range("B68").Value = "CSBBANK"
$losers.Range("C68").Value = -1.4423
$losers.Range("D68").Value = -2.6072
$losers.Range("E68").Value = 1.6693
$losers.Range("B69").Value = "BALAJITELE"
$losers.Range("C69").Value = -1.4302
$losers.Range("D69").Value = -3.6566
$losers.Range("E69").Value = 10.8093
$losers.Range("B70").Value = "JAYNECOIND"
$losers.Range("C70").Value = -1.4297
$losers.Range("D70").Value = 1.8707
$losers.Range("E70").Value = -3.2694
$losers.Range("B71").Value = "BLSE"
$losers.Range("C71").Value = -1.4295
$losers.Range("D71").Value = -0.2687
$losers.Range("E71").Value = 16.8757
$losers.Range("B72").Value = "SUNDROP"
$losers.Range("C72").Value = -1.4214
$losers.Range("D72").Value = -3.2528
$losers.Range("E72").Value = -3.7478
$losers.Range("B73").Value = "WINDMACHIN"
$losers.Range("C73").Value = -1.4157
$losers.Range("D73").Value = -0.3104
$losers.Range("E73").Value = -1.3652
$losers.Range("B74").Value = "INFY"
$losers.Range("C74").Value = -1.4139
$losers.Range("D74").Value = 2.9252
$losers.Range("E74").Value = 2.5433
$losers.Range("B75").Value = "MEGASOFT"
$losers.Range("C75").Value = -1.4042
$losers.Range("D75").Value = -5.7861
$losers.Range("E75").Value = -6.3121
$losers.Range("B76").Value = "EMAMILTD"
$losers.Range("C76").Value = -1.3884
$losers.Range("D76").Value = 0.1653
$losers.Range("E76").Value = -3.6396
